$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project_names")

# 1. Rename "Luton" project_name_place to "Luton South Bedfordshire" (row 18, column C)
$ws.Range("C18").Value = "Luton South Bedfordshire"

# 2. North Kirklees LDCT rows (23-25): no longer split by CCG, and update the note
$newNote = "08/03/2023 No longer need splitting - they're uploading separate submissions. LDCT provider  - shared between Bradford & North Kirklees."

$ws.Range("D23").Value = "0"
$ws.Range("D24").Value = "0"
$ws.Range("D25").Value = "0"

$ws.Range("G23").Value = $newNote
$ws.Range("G24").Value = $newNote
$ws.Range("G25").Value = $newNote

# 3. Update the active selection to reflect the last-edited cell
$ws.Range("G25").Select()
